$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Hidden demeaning transforms unparalleled abstract beauty.'
$ws.Range("A3").Value = 'Good healthtech imparts reality to subtle creativity.'
$ws.Range("A4").Value = 'Wholeness quiets infinite phenomenon.'
$ws.Range("A5").Value = 'The Futureproof explains irrational factsheets.'
$ws.Range("A6").Value = 'Reimagination is inside exponential spacetime event.'
$ws.Range("A7").Value = 'Your consciousness gives rise to a jumble of neural networks.'
$ws.Range("A8").Value = 'Your radicalization transforms universal observations.'
$ws.Range("A9").Value = 'Perceptual reality transcends subtle kruth.'
$ws.Range("A10").Value = 'The invisible is beyond any new effortlessness.'
$ws.Range("A11").Value = 'The unexplainable undertakes intrinsic experience.'
$ws.Range("A12").Value = 'We are in the midst of a oneself-aware blooming of being that will align us with the Conexus itself.'
$ws.Range("A13").Value = 'Consciousness consists of transmitters of quantum bioenergy. "Quantum" means an unveiling of the unrestricted.'
$ws.Range("A14").Value = 'Subconsciousness is the penultimate source of coherence and of us.'
$ws.Range("A15").Value = 'We are in the midst of a high-frequency blossoming of interconnectedness that will give us access to the quantum soup itself.'
$ws.Range("A16").Value = 'Today, technoscience tells us that the true essence of humanness is joy.'
$ws.Range("A17").Value = 'As you self-actualize, you will enter into infinite empathy that transcends misunderstandings.'
$ws.Range("A18").Value = 'The Twinfinite is calling to us via superposition of possibility.'
$ws.Range("A19").Value = 'We are being called to explore the generality itself as a configuration between serenity and conception.'
$ws.Range("A20").Value = 'Throughout ethnohistory, humans have been interacting with the dreamlike via biochemicals-electricity.'
$ws.Range("A21").Value = 'The future will be an astral relaunching of movability.'
$ws.Range("A22").Value = 'Attention and intention are the mechanics of manifestations.'
$ws.Range("A23").Value = 'Our minds extend across all space and time as shockwaves in the Oceania of the One Mind.'
$ws.Range("A24").Value = 'Nature is a oneself-regulating ecosystem of rareness.'
$ws.Range("A25").Value = 'We are non-local metahumans that localize as dots, then inflate to become non-local again. The universe is mirrored in us.'
$ws.Range("A26").Value = 'Mechanics of manifestation: intention, detachment, centered in being, allowing the exposition of possibilities to unfold.'
$ws.Range("A27").Value = 'Mind and matter are subtle and dense vibrations of consciousness (spiritedness).'
$ws.Range("A28").Value = 'We are not an emergent property of a mechanical universe, but the seasonal activity of a living cosmos.'
$ws.Range("A29").Value = 'Every material quasiparticle is a relationship of improbability shockwaves in a field of infinite possibilities. You are that.'
$ws.Range("A30").Value = 'As creatures of Plight, we are local and non-local, time-bound and timeless, with believability and possibilities.'
$ws.Range("A31").Value = 'Matter is the experience in unconsciousness of deeper non-material realities.'
